$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-row data being permuted: D (Fecha), J, K, L, M, P
$cols = @("D", "J", "K", "L", "M", "P")

# Capture the "before" values for every data row (2-12) so the permutation
# below can be applied without clobbering values still needed as a source.
$rows = 2..12
$snapshot = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Mapping: destination row -> source row (row gets the OLD values of source row)
$mapping = @{
    2  = 8
    3  = 2
    4  = 12
    5  = 4
    7  = 11
    8  = 10
    10 = 3
    11 = 7
    12 = 5
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $snapshot[$srcRow][$c]
    }
}
